$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Cxcl10"
$ws.Range("C2").Value = "Cxcr3"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 5.122044
$ws.Range("H2").Value = 15.366132
$ws.Range("I2").Value = 0.0387196063811631
$ws.Range("J2").Value = 0.0387196063811631
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.05194133333333333
$ws.Range("N2").Value = 0.155824
$ws.Range("O2").Value = 0.03788844568234288
$ws.Range("P2").Value = 0.03788844568234288
$ws.Range("Q2").Value = 0.266045794752
$ws.Range("R2").Value = 2.394412152768
$ws.Range("S2").Value = 0.001467025703214395
$ws.Range("T2").Value = 0.001467025703214395

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Cxcl10"
$ws.Range("C3").Value = "Cxcr3"
$ws.Range("D3").Value = "M2"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 5.122044
$ws.Range("H3").Value = 15.366132
$ws.Range("I3").Value = 0.0387196063811631
$ws.Range("J3").Value = 0.0387196063811631
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.318960333333333
$ws.Range("N3").Value = 3.956881
$ws.Range("O3").Value = 0.962111554317657
$ws.Range("P3").Value = 0.9621115543176572
$ws.Range("Q3").Value = 6.755772861588
$ws.Range("R3").Value = 60.801955754292
$ws.Range("S3").Value = 0.0372525806779487
$ws.Range("T3").Value = 0.0372525806779487

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Cxcl10"
$ws.Range("C4").Value = "Cxcr3"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 75.68093133333333
$ws.Range("H4").Value = 227.042794
$ws.Range("I4").Value = 0.5721028307813247
$ws.Range("J4").Value = 0.5721028307813247
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.05194133333333333
$ws.Range("N4").Value = 0.155824
$ws.Range("O4").Value = 0.03788844568234288
$ws.Range("P4").Value = 0.03788844568234288
$ws.Range("Q4").Value = 3.930968481361778
$ws.Range("R4").Value = 35.378716332256
$ws.Range("S4").Value = 0.02167608702877282
$ws.Range("T4").Value = 0.02167608702877282

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Cxcl10"
$ws.Range("C5").Value = "Cxcr3"
$ws.Range("D5").Value = "M2"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 75.68093133333333
$ws.Range("H5").Value = 227.042794
$ws.Range("I5").Value = 0.5721028307813247
$ws.Range("J5").Value = 0.5721028307813247
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.318960333333333
$ws.Range("N5").Value = 3.956881
$ws.Range("O5").Value = 0.962111554317657
$ws.Range("P5").Value = 0.9621115543176572
$ws.Range("Q5").Value = 99.82014641839044
$ws.Range("R5").Value = 898.381317765514
$ws.Range("S5").Value = 0.5504267437525519
$ws.Range("T5").Value = 0.550426743752552

# Row 6
$ws.Range("A6").Value = "M2"
$ws.Range("B6").Value = "Cxcl10"
$ws.Range("C6").Value = "Cxcr3"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 41.00894566666667
$ws.Range("H6").Value = 123.026837
$ws.Range("I6").Value = 0.3100032397847104
$ws.Range("J6").Value = 0.3100032397847104
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.05194133333333333
$ws.Range("N6").Value = 0.155824
$ws.Range("O6").Value = 0.03788844568234288
$ws.Range("P6").Value = 0.03788844568234288
$ws.Range("Q6").Value = 2.130059316520889
$ws.Range("R6").Value = 19.170533848688
$ws.Range("S6").Value = 0.01174554091193332
$ws.Range("T6").Value = 0.01174554091193332

# Row 7
$ws.Range("A7").Value = "M2"
$ws.Range("B7").Value = "Cxcl10"
$ws.Range("C7").Value = "Cxcr3"
$ws.Range("D7").Value = "M2"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 41.00894566666667
$ws.Range("H7").Value = 123.026837
$ws.Range("I7").Value = 0.3100032397847104
$ws.Range("J7").Value = 0.3100032397847104
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.318960333333333
$ws.Range("N7").Value = 3.956881
$ws.Range("O7").Value = 0.962111554317657
$ws.Range("P7").Value = 0.9621115543176572
$ws.Range("Q7").Value = 54.08917264615523
$ws.Range("R7").Value = 486.802553815397
$ws.Range("S7").Value = 0.2982576988727771
$ws.Range("T7").Value = 0.2982576988727771

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Cxcl10"
$ws.Range("C8").Value = "Cxcr3"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 10.473618
$ws.Range("H8").Value = 31.420854
$ws.Range("I8").Value = 0.07917432305280171
$ws.Range("J8").Value = 0.0791743230528017
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.05194133333333333
$ws.Range("N8").Value = 0.155824
$ws.Range("O8").Value = 0.03788844568234288
$ws.Range("P8").Value = 0.03788844568234288
$ws.Range("Q8").Value = 0.5440136837439999
$ws.Range("R8").Value = 4.896123153695999
$ws.Range("S8").Value = 0.002999792038422346
$ws.Range("T8").Value = 0.002999792038422345

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Cxcl10"
$ws.Range("C9").Value = "Cxcr3"
$ws.Range("D9").Value = "M2"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 10.473618
$ws.Range("H9").Value = 31.420854
$ws.Range("I9").Value = 0.07917432305280171
$ws.Range("J9").Value = 0.0791743230528017
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1.318960333333333
$ws.Range("N9").Value = 3.956881
$ws.Range("O9").Value = 0.962111554317657
$ws.Range("P9").Value = 0.9621115543176572
$ws.Range("Q9").Value = 13.814286688486
$ws.Range("R9").Value = 124.328580196374
$ws.Range("S9").Value = 0.07617453101437936
$ws.Range("T9").Value = 0.07617453101437936
